$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110, shifting existing rows 110-117 down to 111-118
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new data record
$ws.Cells.Item(110, 1).Value = 7
$ws.Cells.Item(110, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(110, 3).Value = "Ñuble"
$ws.Cells.Item(110, 4).Value = 44461
$ws.Cells.Item(110, 5).Value = 16
$ws.Cells.Item(110, 6).Value = 100112017
$ws.Cells.Item(110, 7).Value = "Apio"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 160
$ws.Cells.Item(110, 11).Value = 8500
$ws.Cells.Item(110, 12).Value = 9000
$ws.Cells.Item(110, 13).Value = 8750
$ws.Cells.Item(110, 14).Value = "$/docena de matas"
$ws.Cells.Item(110, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(110, 16).Value = 1458
$ws.Cells.Item(110, 17).Value = 6
$ws.Cells.Item(110, 18).Value = "Hortaliza"
